$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows above row 1155, pushing the existing 1155:1176 block
# down to 1159:1180 (matches the new weekly data at the top of the block).
$ws.Rows("1155:1158").Insert()

# New weekly rows: Mercado Mayorista Lo Valledor de Santiago / Melon / Tuna
# dated 2022-03-08 (serial 44628), origin Region Metropolitana.
$newRows = @(
    @{ Row = 1155; I = "Extra";   J = 3200; K = 1200; L = 1300; M = 1241 },
    @{ Row = 1156; I = "Primera"; J = 3800; K = 900;  L = 1000; M = 942  },
    @{ Row = 1157; I = "Segunda"; J = 2300; K = 600;  L = 700;  M = 635  },
    @{ Row = 1158; I = "Tercera"; J = 1700; K = 300;  L = 400;  M = 335  }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 6
    $ws.Cells.Item($row, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
    $ws.Cells.Item($row, 3).Value = "Metropolitana"
    $ws.Cells.Item($row, 4).Value = 44628
    $ws.Cells.Item($row, 5).Value = 13
    $ws.Cells.Item($row, 6).Value = 100112027
    $ws.Cells.Item($row, 7).Value = "Melón"
    $ws.Cells.Item($row, 8).Value = "Tuna"
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = "$/unidad"
    $ws.Cells.Item($row, 15).Value = "Región Metropolitana"
    $ws.Cells.Item($row, 16).Value = $r.M
    $ws.Cells.Item($row, 17).Value = 1
    $ws.Cells.Item($row, 18).Value = "Hortaliza"
}
